# Generate Report for Archive
# The file "547613be-f236-4205-81d9-d6225ab7b667.md" moved from
# "Ready for handoff" to "In Translation" status. Update its Status
# cells on the Overview sheet (zh-cn and de-de columns) as well as on
# the per-locale (zh-cn / de-de) detail sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B4").Value = "In Translation"
$overview.Range("C4").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C4").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C4").Value = "In Translation"
